$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits at the end of the
#    "Nancy for go audit'ing" paragraph (it is hidden from Bookmarks.Count /
#    enumeration, like real Word, but is still addressable by name).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Replace the final "Add more here as need be" list paragraph with:
#      - a new paragraph describing the Apps for Ansible / awscli use, plus
#        the sudo/pip3 install commands (bold), each line separated by a
#        manual line break
#      - a new "Add more here as need b" + <<bookmark>> + "e" paragraph,
#        i.e. the bookmark is re-created inside the final paragraph, split
#        around the last character.
$target = $d.Content
$found = $target.Find.Execute("Add more here as need be", $true, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Add more here as need be' paragraph"
}

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Apps for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ansible</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> use with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>awscli</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> commands:</w:t></w:r><w:r><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> apt install python3-pip</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">pip3 install boto3 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>botocore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t>Add more here as need b</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>e</w:t></w:r></w:p>
'@

$null = $target.InsertXML($newXml)
